$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New leaderboard entries for cycle 17/18 (rows 129-138) ---
# Row 129
$ws.Range("A129").Value = 'NotLogic'
$ws.Range("B129").Value = 17
$ws.Range("C129").Value = 'https://youtu.be/EnUm3o_kcDU?si=czVzhvOhmsCjIzKB'
$ws.Range("D129").Value = 0.1875
$ws.Range("E129").Value = 'Klee'
$ws.Range("F129").Value = 'Bennett'
$ws.Range("G129").Value = 'Xiangling'
$ws.Range("H129").Value = 'Kazuha'
$ws.Range("I129").Value = 'HuTao'
$ws.Range("J129").Value = 'Xingqiu'
$ws.Range("K129").Value = 'Yelan'
$ws.Range("L129").Value = 'Mona'

# Row 130
$ws.Range("A130").Value = 'Yangi'
$ws.Range("B130").Value = 17
$ws.Range("C130").Value = 'https://www.youtube.com/watch?v=LZKeeRT80_E&ab_channel=yangi'
$ws.Range("D130").Value = 0.09722222222222222
$ws.Range("E130").Value = 'Dehya'
$ws.Range("F130").Value = 'Bennett'
$ws.Range("G130").Value = 'Venti'
$ws.Range("H130").Value = 'Rosaria'
$ws.Range("I130").Value = 'Neuvillette'
$ws.Range("J130").Value = 'Xiangling'
$ws.Range("K130").Value = 'Sucrose'
$ws.Range("L130").Value = 'Zhongli'

# Row 131
$ws.Range("A131").Value = 'NotLogic'
$ws.Range("B131").Value = 17
$ws.Range("C131").Value = 'https://youtu.be/HC28jzR18IQ'
$ws.Range("D131").Value = 0.14583333333333334
$ws.Range("E131").Value = 'Klee'
$ws.Range("F131").Value = 'Bennett'
$ws.Range("G131").Value = 'Xiangling'
$ws.Range("H131").Value = 'Kazuha'
$ws.Range("I131").Value = 'Neuvillette'
$ws.Range("J131").Value = 'Kuki'
$ws.Range("K131").Value = 'Collei'
$ws.Range("L131").Value = 'Zhongli'

# Row 132
$ws.Range("A132").Value = 'Dank'
$ws.Range("B132").Value = 18
$ws.Range("C132").Value = 'https://youtu.be/FX-BVi8WNuI '
$ws.Range("D132").Value = 0.25277777777777777
$ws.Range("E132").Value = 'Yoimiya'
$ws.Range("F132").Value = 'Bennett'
$ws.Range("G132").Value = 'Xiangling'
$ws.Range("H132").Value = 'Kazuha'
$ws.Range("I132").Value = 'Raiden'
$ws.Range("J132").Value = 'Yaoyao'
$ws.Range("K132").Value = 'Xingqiu'
$ws.Range("L132").Value = 'Yelan'

# Row 133
$ws.Range("A133").Value = 'NotLogic'
$ws.Range("B133").Value = 18
$ws.Range("C133").Value = 'https://youtu.be/JXQoacbiObk?si=ptzdGCwo0cK4YVs9 '
$ws.Range("D133").Value = 0.22291666666666665
$ws.Range("E133").Value = 'Klee'
$ws.Range("F133").Value = 'Bennett'
$ws.Range("G133").Value = 'Xiangling'
$ws.Range("H133").Value = 'Kazuha'
$ws.Range("I133").Value = 'Neuvillette'
$ws.Range("J133").Value = 'Nahida'
$ws.Range("K133").Value = 'Raiden'
$ws.Range("L133").Value = 'Zhongli'

# Row 134
$ws.Range("A134").Value = 'Tmti'
$ws.Range("B134").Value = 18
$ws.Range("C134").Value = 'https://www.youtube.com/watch?v=F1nd3Bys_h0&t'
$ws.Range("D134").Value = 0.2041666666666667
$ws.Range("E134").Value = 'Yoimiya'
$ws.Range("F134").Value = 'Bennett'
$ws.Range("G134").Value = 'Xiangling'
$ws.Range("H134").Value = 'Kazuha'
$ws.Range("I134").Value = 'HuTao'
$ws.Range("J134").Value = 'Xingqiu'
$ws.Range("K134").Value = 'Yelan'
$ws.Range("L134").Value = 'Jean'

# Row 135
$ws.Range("A135").Value = 'Dank'
$ws.Range("B135").Value = 18
$ws.Range("C135").Value = 'https://youtu.be/kR6vUBqqxR0'
$ws.Range("D135").Value = 0.17569444444444446
$ws.Range("E135").Value = 'HuTao'
$ws.Range("F135").Value = 'Xingqiu'
$ws.Range("G135").Value = 'Sucrose'
$ws.Range("H135").Value = 'Amber'
$ws.Range("I135").Value = 'Neuvillette'
$ws.Range("J135").Value = 'TravelerGeo'
$ws.Range("K135").Value = 'Kazuha'
$ws.Range("L135").Value = 'Diona'

# Row 136
$ws.Range("A136").Value = 'Ghosted'
$ws.Range("B136").Value = 18
$ws.Range("C136").Value = 'https://youtu.be/mNowl4-K3Rg'
$ws.Range("D136").Value = 0.16805555555555554
$ws.Range("E136").Value = 'Ayato'
$ws.Range("F136").Value = 'Rosaria'
$ws.Range("G136").Value = 'Kaeya'
$ws.Range("H136").Value = 'Jean'
$ws.Range("I136").Value = 'Neuvillette'
$ws.Range("J136").Value = 'Fischl'
$ws.Range("K136").Value = 'Beidou'
$ws.Range("L136").Value = 'Sucrose'

# Row 137
$ws.Range("A137").Value = 'Echidna'
$ws.Range("B137").Value = 18
$ws.Range("C137").Value = 'https://youtu.be/Cab8A3zukNc?si=f7HnfISa55Jlg4eq'
$ws.Range("D137").Value = 0.20902777777777778
$ws.Range("E137").Value = 'Lyney'
$ws.Range("F137").Value = 'Bennett'
$ws.Range("G137").Value = 'Albedo'
$ws.Range("H137").Value = 'Zhongli'
$ws.Range("I137").Value = 'Neuvillette'
$ws.Range("J137").Value = 'Kuki'
$ws.Range("K137").Value = 'Collei'
$ws.Range("L137").Value = 'Kirara'

# Row 138
$ws.Range("A138").Value = 'NotLogic'
$ws.Range("B138").Value = 18
$ws.Range("C138").Value = 'https://youtu.be/P9si_fxE0c0?si=MkdXwDLdAylllSU7'
$ws.Range("D138").Value = 0.08611111111111112
$ws.Range("E138").Value = 'HuTao'
$ws.Range("F138").Value = 'Xingqiu'
$ws.Range("G138").Value = 'Sucrose'
$ws.Range("H138").Value = 'Amber'
$ws.Range("I138").Value = 'Neuvillette'
$ws.Range("J138").Value = 'Kuki'
$ws.Range("K138").Value = 'Collei'
$ws.Range("L138").Value = 'Zhongli'

# --- Re-apply the standard row formatting (font/border/alignment/number format) ---
# Rows 127-128 previously used a one-off fill/border variant; bring them back in line
# with the rest of the table, and format the freshly added rows the same way.
$ws.Range("A126:L126").Copy()
$ws.Range("A127:L128").PasteSpecial(-4122)
$ws.Range("A126:L126").Copy()
$ws.Range("A129:L138").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Match the saved selection state from the edit ---
$ws.Range("M138").Select()
